$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values A2:A5 from 5 to 9
$ws.Range("A2").Value = 9
$ws.Range("A3").Value = 9
$ws.Range("A4").Value = 9
$ws.Range("A5").Value = 9

# Update the window size of the workbook view
$win = $excel.ActiveWindow
$win.Width = 30240
$win.Height = 11500

# Update the selection on the sheet to C12
$ws.Range("C12").Select()
